$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.980.81"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.156.51"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.47"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.84"
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.157.04"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  +13.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.17"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.682.07"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.101.84"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.178.06"
$ws.Range("E17").Value = "  +3.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.11"
$ws.Range("E18").Value = "  +4.68%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.36"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.82"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.32"
$ws.Range("E23").Value = "  +4.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.75"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.50"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.91"
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.88"
$ws.Range("E28").Value = "  +8.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.72"
$ws.Range("E30").Value = "  +4.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.76"
$ws.Range("E31").Value = "  +8.78%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.25"
$ws.Range("E34").Value = "  +9.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.54"
$ws.Range("E35").Value = "  +4.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.39"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0896"
$ws.Range("E37").Value = "  +9.62%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0425"
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "464.64"
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  +6.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.67"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.068.23"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  +6.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.70"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0587"
$ws.Range("E47").Value = "  +13.02%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  +5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.89"
$ws.Range("E51").Value = "  +1.31%  "
